$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells that hold numeric-looking or date-looking text must be forced to
# the Text number format first, otherwise Excel will auto-convert them
# into real numbers / dates instead of keeping them as shared strings.
$textCells = @("A2","F2","A3","F3","A4","F4")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 (record for Shane Paz) -- TEST MARKER
$ws.Range("A2").Value = "183"
$ws.Range("D2").Value = "Shane"
$ws.Range("E2").Value = "Paz"
$ws.Range("F2").Value = "1972-06-03"
$ws.Range("G2").Value = "471 Waterview LanePena Blanca, NM 87041"
$ws.Range("H2").Value = "505-465-5925"
$ws.Range("J2").Value = "Simply Appraisals"
$ws.Range("K2").Value = "Orthotics technician"
$ws.Range("L2").Value = "sieKa5uDie"

# Row 3 (record for Sylvia Allen)
$ws.Range("A3").Value = "174"
$ws.Range("C3").Value = "Australia"
$ws.Range("D3").Value = "Sylvia"
$ws.Range("E3").Value = "Allen"
$ws.Range("F3").Value = "2003-05-24"
$ws.Range("G3").Value = "54 Wallis StreetROSE BAY NORTH NSW 2030"
$ws.Range("H3").Value = "(02) 9579 8213"
$ws.Range("J3").Value = "Cut Rite"
$ws.Range("K3").Value = "Webmaster"
$ws.Range("L3").Value = "ahxeiJoo6"

# Row 4 (record for Carol Schmidt)
$ws.Range("A4").Value = "155"
$ws.Range("C4").Value = "New Zealands"
$ws.Range("D4").Value = "Carol"
$ws.Range("E4").Value = "Schmidt"
$ws.Range("F4").Value = "1960-02-02"
$ws.Range("G4").Value = "7 Hautana StreetBoulcottLower Hutt 5010"
$ws.Range("H4").Value = "(028) 6342-658"
$ws.Range("J4").Value = "Incluesiv"
$ws.Range("K4").Value = "Museum director"
$ws.Range("L4").Value = "aiNguLaiN9"
